# Auto-generated Excel COM-interop script to apply the Titan_Profits market-data refresh.
# Updates cached price/profit columns (H:N) for specific Leve rows across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Cells.Item(34, 8).Value = 10543
$ws.Cells.Item(34, 9).Value = 3477.7144
$ws.Cells.Item(34, 10).Value = 60000
$ws.Cells.Item(34, 11).Value = 3477.7144
$ws.Cells.Item(34, 12).Value = 60000
$ws.Cells.Item(34, 13).Value = -3274.7144
$ws.Cells.Item(34, 14).Value = -60406

# Row 36
$ws.Cells.Item(36, 8).Value = 10543
$ws.Cells.Item(36, 9).Value = 3477.7144
$ws.Cells.Item(36, 10).Value = 60000
$ws.Cells.Item(36, 11).Value = 3477.7144
$ws.Cells.Item(36, 12).Value = 60000
$ws.Cells.Item(36, 13).Value = -2762.7144
$ws.Cells.Item(36, 14).Value = -61430

# Row 132
$ws.Cells.Item(132, 8).Value = 307427.47
$ws.Cells.Item(132, 9).Value = 348274.1
$ws.Cells.Item(132, 10).Value = 21501.2
$ws.Cells.Item(132, 11).Value = 1044822.3
$ws.Cells.Item(132, 12).Value = 64503.60000000001
$ws.Cells.Item(132, 13).Value = -1042292.3
$ws.Cells.Item(132, 14).Value = -69563.60000000001

# Row 138
$ws.Cells.Item(138, 8).Value = 1181.99
$ws.Cells.Item(138, 9).Value = 572.4828
$ws.Cells.Item(138, 10).Value = 2023.6904
$ws.Cells.Item(138, 11).Value = 1717.4484
$ws.Cells.Item(138, 12).Value = 6071.0712
$ws.Cells.Item(138, 13).Value = 3422.5516
$ws.Cells.Item(138, 14).Value = -16351.0712

# Row 141
$ws.Cells.Item(141, 8).Value = 2076.932
$ws.Cells.Item(141, 9).Value = 1302.2463
$ws.Cells.Item(141, 10).Value = 4890.263
$ws.Cells.Item(141, 11).Value = 3906.7389
$ws.Cells.Item(141, 12).Value = 14670.789
$ws.Cells.Item(141, 13).Value = 1273.2611
$ws.Cells.Item(141, 14).Value = -25030.789

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 17366.508
$ws.Cells.Item(32, 9).Value = 2711.4922
$ws.Cells.Item(32, 10).Value = 255510.5
$ws.Cells.Item(32, 11).Value = 2711.4922
$ws.Cells.Item(32, 12).Value = 255510.5
$ws.Cells.Item(32, 13).Value = -2424.4922

# Row 74
$ws.Cells.Item(74, 8).Value = 3409.7932
$ws.Cells.Item(74, 9).Value = 1143.0233
$ws.Cells.Item(74, 10).Value = 9907.866
$ws.Cells.Item(74, 11).Value = 1143.0233
$ws.Cells.Item(74, 12).Value = 9907.866
$ws.Cells.Item(74, 13).Value = -269.0233000000001
$ws.Cells.Item(74, 14).Value = -11655.866

# Row 77
$ws.Cells.Item(77, 8).Value = 3409.7932
$ws.Cells.Item(77, 9).Value = 1143.0233
$ws.Cells.Item(77, 10).Value = 9907.866
$ws.Cells.Item(77, 11).Value = 5715.1165
$ws.Cells.Item(77, 12).Value = 49539.33
$ws.Cells.Item(77, 13).Value = -1347.1165
$ws.Cells.Item(77, 14).Value = -58275.33

# Row 132
$ws.Cells.Item(132, 8).Value = 3073.861
$ws.Cells.Item(132, 9).Value = 2818
$ws.Cells.Item(132, 10).Value = 3841.4443
$ws.Cells.Item(132, 11).Value = 8454
$ws.Cells.Item(132, 12).Value = 11524.3329
$ws.Cells.Item(132, 13).Value = -5924

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 1502.6
$ws.Cells.Item(107, 9).Value = 833.3333
$ws.Cells.Item(107, 10).Value = 2506.5
$ws.Cells.Item(107, 11).Value = 833.3333
$ws.Cells.Item(107, 12).Value = 2506.5
$ws.Cells.Item(107, 13).Value = 1086.6667
$ws.Cells.Item(107, 14).Value = -6346.5

# Row 119
$ws.Cells.Item(119, 8).Value = 41230.5
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 41230.5
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 41230.5
$ws.Cells.Item(119, 14).Value = -50906.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 23617.318
$ws.Cells.Item(16, 9).Value = 34119.2
$ws.Cells.Item(16, 10).Value = 1113.2858
$ws.Cells.Item(16, 11).Value = 34119.2
$ws.Cells.Item(16, 12).Value = 1113.2858
$ws.Cells.Item(16, 13).Value = -33832.2
$ws.Cells.Item(16, 14).Value = -1687.2858

# Row 105
$ws.Cells.Item(105, 8).Value = 916.5
$ws.Cells.Item(105, 9).Value = 900
$ws.Cells.Item(105, 10).Value = 999
$ws.Cells.Item(105, 11).Value = 900
$ws.Cells.Item(105, 12).Value = 999
$ws.Cells.Item(105, 13).Value = 847
$ws.Cells.Item(105, 14).Value = -4493

# Row 113
$ws.Cells.Item(113, 8).Value = 23617.318
$ws.Cells.Item(113, 9).Value = 34119.2
$ws.Cells.Item(113, 10).Value = 1113.2858
$ws.Cells.Item(113, 11).Value = 34119.2
$ws.Cells.Item(113, 12).Value = 1113.2858
$ws.Cells.Item(113, 13).Value = -31949.2
$ws.Cells.Item(113, 14).Value = -5453.2858

# Row 123
$ws.Cells.Item(123, 8).Value = 31000
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 31000
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 31000
$ws.Cells.Item(123, 14).Value = -40800

# Row 132
$ws.Cells.Item(132, 8).Value = 2354.1794
$ws.Cells.Item(132, 9).Value = 1921.2424
$ws.Cells.Item(132, 10).Value = 4735.3335
$ws.Cells.Item(132, 11).Value = 5763.7272
$ws.Cells.Item(132, 12).Value = 14206.0005
$ws.Cells.Item(132, 13).Value = -3233.7272
$ws.Cells.Item(132, 14).Value = -19266.0005

# Row 134
$ws.Cells.Item(134, 8).Value = 1883.9305
$ws.Cells.Item(134, 9).Value = 1231
$ws.Cells.Item(134, 10).Value = 4588.9287
$ws.Cells.Item(134, 11).Value = 3693
$ws.Cells.Item(134, 12).Value = 13766.7861
$ws.Cells.Item(134, 13).Value = -1158
$ws.Cells.Item(134, 14).Value = -18836.7861

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Cells.Item(23, 8).Value = 773.1053000000001
$ws.Cells.Item(23, 9).Value = 3500.3333
$ws.Cells.Item(23, 10).Value = 261.75
$ws.Cells.Item(23, 11).Value = 10500.9999
$ws.Cells.Item(23, 12).Value = 785.25
$ws.Cells.Item(23, 13).Value = -10265.9999
$ws.Cells.Item(23, 14).Value = -1255.25

# Row 75
$ws.Cells.Item(75, 8).Value = 1000
$ws.Cells.Item(75, 9).Value = 1000
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 3000
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = -2002
$ws.Cells.Item(75, 14).ClearContents()

# Row 78
$ws.Cells.Item(78, 8).Value = 1000
$ws.Cells.Item(78, 9).Value = 1000
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 9000
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 13).Value = -4008
$ws.Cells.Item(78, 14).ClearContents()

# Row 80
$ws.Cells.Item(80, 8).Value = 1198.5714
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 1198.5714
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 3595.7142
$ws.Cells.Item(80, 14).Value = -5467.7142

# Row 82
$ws.Cells.Item(82, 8).Value = 1881.8182
$ws.Cells.Item(82, 9).Value = 200
$ws.Cells.Item(82, 10).Value = 2050
$ws.Cells.Item(82, 11).Value = 600
$ws.Cells.Item(82, 12).Value = 6150
$ws.Cells.Item(82, 13).Value = -194
$ws.Cells.Item(82, 14).Value = -6962

# Row 83
$ws.Cells.Item(83, 8).Value = 1198.5714
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 1198.5714
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 10787.1426
$ws.Cells.Item(83, 14).Value = -20147.1426

# Row 85
$ws.Cells.Item(85, 8).Value = 1881.8182
$ws.Cells.Item(85, 9).Value = 200
$ws.Cells.Item(85, 10).Value = 2050
$ws.Cells.Item(85, 11).Value = 600
$ws.Cells.Item(85, 12).Value = 6150
$ws.Cells.Item(85, 13).Value = 804
$ws.Cells.Item(85, 14).Value = -8958

# Row 87
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).ClearContents()

# Row 88
$ws.Cells.Item(88, 8).Value = 5000
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 5000
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 15000
$ws.Cells.Item(88, 14).Value = -15856

# Row 90
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 13).ClearContents()

# Row 91
$ws.Cells.Item(91, 8).Value = 5000
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 5000
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 15000
$ws.Cells.Item(91, 14).Value = -17964

# Row 107
$ws.Cells.Item(107, 8).Value = 904.4074000000001
$ws.Cells.Item(107, 9).Value = 1211.25
$ws.Cells.Item(107, 10).Value = 458.0909
$ws.Cells.Item(107, 11).Value = 3633.75
$ws.Cells.Item(107, 12).Value = 1374.2727
$ws.Cells.Item(107, 13).Value = -1713.75
$ws.Cells.Item(107, 14).Value = -5214.2727

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Cells.Item(15, 8).Value = 5000
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 5000
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 5000
$ws.Cells.Item(15, 14).Value = -5576

# Row 81
$ws.Cells.Item(81, 8).Value = 5000
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 5000
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 12).Value = 5000
$ws.Cells.Item(81, 14).Value = -6996

# Row 84
$ws.Cells.Item(84, 8).Value = 5000
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 5000
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 15000
$ws.Cells.Item(84, 14).Value = -24984

# Row 122
$ws.Cells.Item(122, 8).Value = 484497.3
$ws.Cells.Item(122, 9).Value = 586091.6
$ws.Cells.Item(122, 10).Value = 1924.25
$ws.Cells.Item(122, 11).Value = 1758274.8
$ws.Cells.Item(122, 12).Value = 5772.75
$ws.Cells.Item(122, 13).Value = -1755824.8
$ws.Cells.Item(122, 14).Value = -10672.75

# Row 132
$ws.Cells.Item(132, 8).Value = 2391.6462
$ws.Cells.Item(132, 9).Value = 2141.82
$ws.Cells.Item(132, 10).Value = 3224.4
$ws.Cells.Item(132, 11).Value = 6425.460000000001
$ws.Cells.Item(132, 12).Value = 9673.200000000001
$ws.Cells.Item(132, 13).Value = -3895.460000000001
$ws.Cells.Item(132, 14).Value = -14733.2

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Cells.Item(132, 8).Value = 6534.294
$ws.Cells.Item(132, 9).Value = 8767.9375
$ws.Cells.Item(132, 10).Value = 4548.8335
$ws.Cells.Item(132, 11).Value = 26303.8125
$ws.Cells.Item(132, 12).Value = 13646.5005
$ws.Cells.Item(132, 13).Value = -23773.8125
$ws.Cells.Item(132, 14).Value = -18706.5005

# Row 136
$ws.Cells.Item(136, 8).Value = 3945.1333
$ws.Cells.Item(136, 9).Value = 2137.9062
$ws.Cells.Item(136, 10).Value = 8393.691999999999
$ws.Cells.Item(136, 11).Value = 6413.7186
$ws.Cells.Item(136, 12).Value = 25181.076
$ws.Cells.Item(136, 13).Value = -3863.7186
$ws.Cells.Item(136, 14).Value = -30281.076

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 2416746
$ws.Cells.Item(107, 9).Value = 5051714
$ws.Cells.Item(107, 10).Value = 1358.5834
$ws.Cells.Item(107, 11).Value = 15155142
$ws.Cells.Item(107, 12).Value = 4075.7502
$ws.Cells.Item(107, 13).Value = -15153222
$ws.Cells.Item(107, 14).Value = -7915.7502

# Row 122
$ws.Cells.Item(122, 8).Value = 72369.64
$ws.Cells.Item(122, 9).Value = 91663.17999999999
$ws.Cells.Item(122, 10).Value = 1626.6666
$ws.Cells.Item(122, 11).Value = 274989.54
$ws.Cells.Item(122, 12).Value = 4879.9998
$ws.Cells.Item(122, 13).Value = -272539.54

# Row 126
$ws.Cells.Item(126, 8).Value = 63683.312
$ws.Cells.Item(126, 9).Value = 77979.46000000001
$ws.Cells.Item(126, 10).Value = 1733.3334
$ws.Cells.Item(126, 11).Value = 233938.38
$ws.Cells.Item(126, 12).Value = 5200.0002
$ws.Cells.Item(126, 13).Value = -231468.38

# Row 132
$ws.Cells.Item(132, 8).Value = 7144629
$ws.Cells.Item(132, 9).Value = 9805496
$ws.Cells.Item(132, 10).Value = 2301
$ws.Cells.Item(132, 11).Value = 29416488
$ws.Cells.Item(132, 12).Value = 6903
$ws.Cells.Item(132, 13).Value = -29413958
$ws.Cells.Item(132, 14).Value = -11963

# Row 136
$ws.Cells.Item(136, 8).Value = 18619.824
$ws.Cells.Item(136, 9).Value = 21772.361
$ws.Cells.Item(136, 10).Value = 3802.9
$ws.Cells.Item(136, 11).Value = 65317.083
$ws.Cells.Item(136, 12).Value = 11408.7
$ws.Cells.Item(136, 13).Value = -62767.083
$ws.Cells.Item(136, 14).Value = -16508.7
